# Apply updated cryptocurrency price/volume data to worksheet.
# A leading apostrophe forces Excel to store these as text values
# (matching the source inlineStr cell type) instead of auto-converting
# numeric-looking strings (e.g. "547.33", "1.00") into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.789.73"
$ws.Range("E2").Value = "'  -4.16%  "
$ws.Range("D3").Value = "'2.445.24"
$ws.Range("E3").Value = "'  -6.05%  "
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'547.33"
$ws.Range("E5").Value = "'  -4.09%  "
$ws.Range("D6").Value = "'144.42"
$ws.Range("E6").Value = "'  -6.10%  "
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "'  -3.47%  "
$ws.Range("D9").Value = "'2.442.68"
$ws.Range("E9").Value = "'  -6.08%  "
$ws.Range("E10").Value = "'  -7.16%  "
$ws.Range("E11").Value = "'  -1.78%  "
$ws.Range("D12").Value = "'5.35"
$ws.Range("E12").Value = "'  -7.49%  "
$ws.Range("D13").Value = "'0.351"
$ws.Range("E13").Value = "'  -6.65%  "
$ws.Range("D14").Value = "'25.89"
$ws.Range("E14").Value = "'  -7.07%  "
$ws.Range("D15").Value = "'2.892.18"
$ws.Range("E15").Value = "'  -5.90%  "
$ws.Range("E16").Value = "'  -8.49%  "
$ws.Range("D17").Value = "'60.714.21"
$ws.Range("D18").Value = "'2.450.15"
$ws.Range("E18").Value = "'  -5.75%  "
$ws.Range("D19").Value = "'11.01"
$ws.Range("E19").Value = "'  -7.27%  "
$ws.Range("D20").Value = "'6.89"
$ws.Range("E20").Value = "'  -7.53%  "
$ws.Range("D21").Value = "'4.15"
$ws.Range("E21").Value = "'  -6.78%  "
$ws.Range("D22").Value = "'317.79"
$ws.Range("E22").Value = "'  -6.29%  "
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("D24").Value = "'63.30"
$ws.Range("E24").Value = "'  -5.53%  "
$ws.Range("E25").Value = "'  -1.24%  "
$ws.Range("B26").Value = "'WrappedeETH"
$ws.Range("C26").Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "'2.573.45"
$ws.Range("E26").Value = "'  -5.69%  "
$ws.Range("B27").Value = "'PEPE"
$ws.Range("C27").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "'0.0₃0960"
$ws.Range("E27").Value = "'  -8.45%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "'  -0.27%  "
$ws.Range("B29").Value = "'Fetch.AI"
$ws.Range("C29").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "'1.48"
$ws.Range("E29").Value = "'  -3.45%  "
$ws.Range("B30").Value = "'Bittensor"
$ws.Range("C30").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "'534.30"
$ws.Range("E30").Value = "'  -7.17%  "
$ws.Range("D31").Value = "'8.30"
$ws.Range("E31").Value = "'  -8.11%  "
$ws.Range("D32").Value = "'7.56"
$ws.Range("E32").Value = "'  -2.40%  "
$ws.Range("D33").Value = "'0.148"
$ws.Range("E33").Value = "'  -7.23%  "
$ws.Range("D34").Value = "'1.88"
$ws.Range("E34").Value = "'  -7.29%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "'  -8.30%  "
$ws.Range("D36").Value = "'5.80"
$ws.Range("E36").Value = "'  -10.61%  "
$ws.Range("B37").Value = "'FirstDigitalUSD"
$ws.Range("C37").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "'  +0.04%  "
$ws.Range("B38").Value = "'NEARProtocol"
$ws.Range("C38").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.82"
$ws.Range("E38").Value = "'  -9.46%  "
$ws.Range("D39").Value = "'0.375"
$ws.Range("E39").Value = "'  -5.79%  "
$ws.Range("D40").Value = "'18.39"
$ws.Range("E40").Value = "'  -5.96%  "
$ws.Range("D41").Value = "'144.72"
$ws.Range("E41").Value = "'  -6.31%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "'  +0.03%  "
$ws.Range("D43").Value = "'1.70"
$ws.Range("E43").Value = "'  -8.06%  "
$ws.Range("D44").Value = "'39.72"
$ws.Range("E44").Value = "'  -4.30%  "
$ws.Range("D45").Value = "'2.28"
$ws.Range("E45").Value = "'  -7.56%  "
$ws.Range("D46").Value = "'145.16"
$ws.Range("E46").Value = "'  -7.49%  "
$ws.Range("D47").Value = "'3.54"
$ws.Range("E47").Value = "'  -7.24%  "
$ws.Range("D48").Value = "'20.67"
$ws.Range("E48").Value = "'  -11.11%  "
$ws.Range("D49").Value = "'0.0526"
$ws.Range("E49").Value = "'  -8.90%  "
$ws.Range("D50").Value = "'0.581"
$ws.Range("E50").Value = "'  -7.15%  "
$ws.Range("D51").Value = "'0.0938"
$ws.Range("E51").Value = "'  -5.34%  "
